$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

# Ordered so that "95÷7=" (old at one cell) is replaced before it is
# reintroduced as a new value elsewhere ("31÷7=" -> "95÷7="); this avoids
# the later insertion being clobbered by an earlier find/replace pass.
Replace-Text "53÷4=" "39÷2="
Replace-Text "92÷5=" "67÷9="
Replace-Text "98÷2=" "68÷4="
Replace-Text "97÷4=" "68÷9="
Replace-Text "84÷2=" "54÷3="
Replace-Text "18÷9=" "57÷2="
Replace-Text "18÷5=" "46÷2="
Replace-Text "15÷3=" "48÷8="
Replace-Text "67÷2=" "31÷2="
Replace-Text "72÷3=" "99÷5="
Replace-Text "95÷7=" "19÷7="
Replace-Text "24÷3=" "37÷3="
Replace-Text "44÷2=" "42÷4="
Replace-Text "49÷7=" "42÷4="
Replace-Text "31÷7=" "95÷7="
Replace-Text "84÷5=" "51÷5="
Replace-Text "28÷6=" "40÷3="
Replace-Text "76÷4=" "69÷2="
Replace-Text "66÷4=" "68÷9="
Replace-Text "93÷4=" "65÷7="
Replace-Text "98÷7=" "34÷4="
Replace-Text "88÷4=" "18÷3="
Replace-Text "22÷4=" "52÷5="
Replace-Text "71÷7=" "61÷8="
Replace-Text "48÷9=" "34÷8="
